# Applies the scrape-refresh edit described by the diff:
#  - A new item (Magic Matic navy ecorepel, previously scraped at row 7) now
#    sorts first and occupies row 3; the four rows that used to sit at rows 3-6
#    (Bustier, T-Shirt, Söckchen, Strumpfhose) each shift down by one row (4-7).
#  - Every row timestamp (column O) is refreshed to the new crawl time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# id/price/priceContextPrice columns look numeric but are stored as plain text
# in this sheet - prefix with a quote (like typing '123 into a cell) so Excel
# keeps them as text instead of converting them to numbers.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
}

# Row 3: Magic Matic navy ecorepel
Set-TextValue $ws.Range("A3") '5799901001'
$ws.Range("B3").Value = 'Magic Matic navy ecorepel'
$ws.Range("C3").Value = '/de/haushalt-tier/bekleidung/taschen-accessoires/schirme/magic-matic-navy-ecorepel/p/5799901001'
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 'Coop'
Set-TextValue $ws.Range("H3") '19.95'
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = '[''haushalt-tier'', ''bekleidung'', ''taschen-accessoires'', ''schirme'']'
$ws.Range("N3").Value = 'Magic Matic navy ecorepel 19.95 Schweizer Franken'

# Row 4: Naturaline Damen Bustier Weiss L
Set-TextValue $ws.Range("A4") '3875554009'
$ws.Range("B4").Value = 'Naturaline Damen Bustier Weiss L'
$ws.Range("C4").Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-damen-bustier-weiss-l/p/3875554009'
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 'Coop'
Set-TextValue $ws.Range("H4") '14.95'
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''unterwaesche'']'
$ws.Range("N4").Value = 'Naturaline Damen Bustier Weiss L 14.95 Schweizer Franken'

# Row 5: Naturaline Herren T-Shirt Kurzarm schwarzXL
Set-TextValue $ws.Range("A5") '6031467019'
$ws.Range("B5").Value = 'Naturaline Herren T-Shirt Kurzarm schwarzXL'
$ws.Range("C5").Value = '/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-schwarzxl/p/6031467019'
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 'Coop'
Set-TextValue $ws.Range("H5") '24.95'
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = '[''haushalt-tier'', ''bekleidung'', ''shirts-pullover'', ''herren-shirt'']'
$ws.Range("N5").Value = 'Naturaline Herren T-Shirt Kurzarm schwarzXL 24.95 Schweizer Franken'

# Row 6: Avela Söckchen Pure Noir One Size
Set-TextValue $ws.Range("A6") '6077154005'
$ws.Range("B6").Value = 'Avela Söckchen Pure Noir One Size'
$ws.Range("C6").Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/socken/avela-soeckchen-pure-noir-one-size/p/6077154005'
$ws.Range("D6").Value = '2ST'
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 'Avela'
Set-TextValue $ws.Range("H6") '2.95'
$ws.Range("I6").Value = '1.48/1ST'
$ws.Range("J6").Value = 'Preis pro 1 Stück'
Set-TextValue $ws.Range("K6") '1.48'
$ws.Range("L6").Value = '1ST'
$ws.Range("M6").Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''socken'']'
$ws.Range("N6").Value = 'Avela Söckchen Pure Noir One Size 2.95 Schweizer Franken'

# Row 7: Avela Strumpfhose Madame Natural  11 - 12
Set-TextValue $ws.Range("A7") '6075745012'
$ws.Range("B7").Value = 'Avela Strumpfhose Madame Natural  11 - 12'
$ws.Range("C7").Value = '/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-natural-11-12/p/6075745012'
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 'Avela'
Set-TextValue $ws.Range("H7") '5.95'
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = '[''haushalt-tier'', ''bekleidung'', ''socken-unterwaesche'', ''struempfe'']'
$ws.Range("N7").Value = 'Avela Strumpfhose Madame Natural  11 - 12 5.95 Schweizer Franken'

# Refresh the scrape timestamp for every data row (2 through 73)
$newTimestamp = '2022-08-21 20:58:15'
for ($r = 2; $r -le 73; $r++) {
    $ws.Range("O$r").Value = $newTimestamp
}
